$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Change header text of C1 from "month" to "year" (سابقه(ماه) -> سابقه(سال))
$ws.Range("C1").Value = "سابقه(سال)"

# 2. Update selection to L9
$ws.Range("L9").Select()

# 3. Replace data validations on column C:
#    - remove the C1-only whole-number validation
#    - remove the C2:C1001 whole-number validation (allowBlank=false)
#    - add a single combined C1:C1001 whole-number validation (allowBlank=true)
$ws.Range("C1").Validation.Delete()
$ws.Range("C2:C1001").Validation.Delete()

$ws.Range("C1:C1001").Validation.Add(1, 1, 1, "0", "1E+034")
$ws.Range("C1:C1001").Validation.IgnoreBlank = $true
$ws.Range("C1:C1001").Validation.InCellDropdown = $true
$ws.Range("C1:C1001").Validation.ShowInput = $true
$ws.Range("C1:C1001").Validation.ShowError = $true
